# FDF format document update
#  - Split the "data de criacao (em base64)" paragraph into three
#    separate paragraphs (criacao / inicio / termino), moving the
#    "_GoBack" bookmark onto the new "termino" paragraph.
#  - Move a few <w:lastRenderedPageBreak/> markers around to reflect the
#    new pagination caused by the extra paragraphs.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Get-ParaIndexByTextAfter($doc, $text, $prevText) {
    for ($i = 2; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            $p = $doc.Paragraphs.Item($i - 1).Range.Text
            $p = $p.TrimEnd([char]13, [char]7)
            if ($p -eq $prevText) {
                return $i
            }
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) "data de criacao (em base64)" -> three paragraphs
# ---------------------------------------------------------------------

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'
$pPr = '<w:pPr><w:ind w:left="708"/>' + $rPr + '</w:pPr>'

$newXml = '<w:p ' + $wNs + '>' + $pPr + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>data</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> de </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>criacao</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> (em </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>base64</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>)</w:t></w:r>' + `
    '</w:p>' + `
'<w:p ' + $wNs + '>' + $pPr + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>data</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> de início do formulário (em base64)</w:t></w:r>' + `
    '</w:p>' + `
'<w:p ' + $wNs + '>' + $pPr + `
    '<w:r>' + $rPr + '<w:t>data de término do formulário (em base64)</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$idx = Get-ParaIndexByText $d "data de criacao (em base64)"
if ($idx -eq -1) {
    throw "Could not locate the 'data de criacao (em base64)' paragraph"
}
$d.Paragraphs.Item($idx).Range.InsertXML($newXml)

# ---------------------------------------------------------------------
# 2) move the <w:lastRenderedPageBreak/> markers
# ---------------------------------------------------------------------

# "nome" (the field description, not the "nome=..." example further down)
# gains a page-break marker in front of its run now that the header grew
# by two paragraphs.
$idx = Get-ParaIndexByText $d "nome=nome do formulário"
if ($idx -eq -1) { throw "Could not locate 'nome=nome do formulário'" }
$xml = '<w:p ' + $wNs + '><w:pPr><w:ind w:left="708"/>' + $rPr + '</w:pPr>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>nome</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>=nome do formulário</w:t></w:r>' + `
    '</w:p>'
$d.Paragraphs.Item($idx).Range.InsertXML($xml)

# "local_respostas=local aonde..." loses its page-break marker.
$idx = Get-ParaIndexByText $d "local_respostas=local aonde está localizado o campo de respostas"
if ($idx -eq -1) { throw "Could not locate 'local_respostas=local aonde...'" }
$xml = '<w:p ' + $wNs + '><w:pPr><w:ind w:left="708"/>' + $rPr + '</w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>local</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>_respostas</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>=local aonde está localizado o campo de respostas</w:t></w:r>' + `
    '</w:p>'
$d.Paragraphs.Item($idx).Range.InsertXML($xml)

# "= tipo,id" right after "Se tipo for EXCLUSIVA (pode escolher apenas uma)"
# loses its page-break marker.
$idx = Get-ParaIndexByTextAfter $d "= tipo,id" "Se tipo for EXCLUSIVA (pode escolher apenas uma)"
if ($idx -eq -1) { throw "Could not locate '= tipo,id' after EXCLUSIVA" }
$xml = '<w:p ' + $wNs + '><w:pPr><w:ind w:left="708"/>' + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">= </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t>tipo,id</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$d.Paragraphs.Item($idx).Range.InsertXML($xml)

# "- 35" (example file, ALTERNATIVA answers) gains a page-break marker.
$rPr24 = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$idx = Get-ParaIndexByText $d "- 35"
if ($idx -eq -1) { throw "Could not locate '- 35'" }
$xml = '<w:p ' + $wNs + '><w:pPr>' + $rPr24 + '</w:pPr>' + `
    '<w:r>' + $rPr24 + '<w:lastRenderedPageBreak/><w:t>- 35</w:t></w:r>' + `
    '</w:p>'
$d.Paragraphs.Item($idx).Range.InsertXML($xml)

# "=EXCLUSIVA,30" (example file) loses its page-break marker.
$idx = Get-ParaIndexByText $d "=EXCLUSIVA,30"
if ($idx -eq -1) { throw "Could not locate '=EXCLUSIVA,30'" }
$xml = '<w:p ' + $wNs + '><w:pPr>' + $rPr24 + '</w:pPr>' + `
    '<w:r>' + $rPr24 + '<w:t>=EXCLUSIVA,30</w:t></w:r>' + `
    '</w:p>'
$d.Paragraphs.Item($idx).Range.InsertXML($xml)

Write-Host "Done."
